$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last refreshed" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Junio de 2020 a las 04:50"

# --- Countries around Irak/Barein/Japon/Austria/Bolivia (rows 50-53) ---
# Bolivia's stats got refreshed and the row now sits right after Irak,
# pushing Barein/Japon/Austria down one row each.
$ws.Range("A50").Value = "Bolivia"
$ws.Range("B50").Value = 17842
$ws.Range("C50").Value = 913
$ws.Range("D50").Value = 2768
$ws.Range("E50").Value = 14489
$ws.Range("G50").Value = 26
$ws.Range("H50").Value = 585

$ws.Range("A51").Value = "Barein"
$ws.Range("B51").Value = 17713
$ws.Range("D51").Value = 12191
$ws.Range("E51").Value = 5485
$ws.Range("H51").Value = 37

$ws.Range("A52").Value = "Japon"
$ws.Range("B52").Value = 17382
$ws.Range("D52").Value = 15580
$ws.Range("E52").Value = 878
$ws.Range("H52").Value = 924

$ws.Range("A53").Value = "Austria"
$ws.Range("B53").Value = 17078
$ws.Range("D53").Value = 16012
$ws.Range("E53").Value = 389
$ws.Range("H53").Value = 677

# --- Camerun stats refresh (row 71), no reordering ---
$ws.Range("B71").Value = 7320
$ws.Range("C71").Value = 18
$ws.Range("D71").Value = 6815
$ws.Range("E71").Value = 403

# --- Countries around Yibuti/Hungria/Luxemburgo/Haiti (rows 82-84) ---
# Haiti's stats got refreshed and the row now sits right after Republica
# de Yibuti, pushing Hungria/Luxemburgo down one row each.
$ws.Range("A82").Value = "Haiti"
$ws.Range("B82").Value = 4165
$ws.Range("C82").Value = 224
$ws.Range("D82").Value = 24
$ws.Range("E82").Value = 4071
$ws.Range("G82").Value = 6
$ws.Range("H82").Value = 70

$ws.Range("A83").Value = "Hungria"
$ws.Range("B83").Value = 4064
$ws.Range("D83").Value = 2476
$ws.Range("E83").Value = 1029
$ws.Range("H83").Value = 559

$ws.Range("A84").Value = "Luxemburgo"
$ws.Range("B84").Value = 4063
$ws.Range("D84").Value = 3922
$ws.Range("E84").Value = 31
$ws.Range("H84").Value = 110

# --- Sahara Occidental stats refresh (row 212), no reordering ---
$ws.Range("D212").Value = 8
$ws.Range("E212").Value = 0

# --- Swap Papua Nueva Guinea / Islas Virgenes Britanicas (rows 213-214) ---
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
